$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample")

# Origin column (E) changes from "HIR" to "BNE" for the four fare rows
$ws.Range("E2").Value = "BNE"
$ws.Range("E3").Value = "BNE"
$ws.Range("E4").Value = "BNE"
$ws.Range("E5").Value = "BNE"

# Update the active selection to match the edited workbook's saved cursor position
$ws.Range("H4").Select()
